$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header for column E
$ws.Range("E1").Value = "Address"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null

# Set column E width
$ws.Columns.Item(5).ColumnWidth = 97.65

# Fill address data for each park row (rows 2-44)
$ws.Range("E2").Value = "Adventure Island, Western Esplanade, Southend-on-Sea SS1 1EE, United Kingdom"
$ws.Range("E3").Value = "Adventure Wonderland, Merritown Ln, Hurn, Christchurch BH23 6BA"
$ws.Range("E4").Value = "Barry Island Pleasure Park, Friars Rd, Barry CF62 5TR"
$ws.Range("E5").Value = "Blackgang Chine, Ventnor PO38 2HN"
$ws.Range("E6").Value = "Brean Leisure Park, Richard's Wy, Brean, Burnham-on-Sea TA8 2RA, United Kingdom"
$ws.Range("E7").Value = "Camel Creek Family Adventure Park, Tredinnick, Wadebridge PL27 7RA"
$ws.Range("E8").Value = "Address not found"
$ws.Range("E9").Value = "Address not found"
$ws.Range("E10").Value = "Address not found"
$ws.Range("E11").Value = "Dreamland Margate, 49-51, Marine Terrace, Margate CT9 1XJ"
$ws.Range("E12").Value = "Address not found"
$ws.Range("E13").Value = "Fantasy Island, Sea Ln, Ingoldmells, Skegness PE25 1RH, United Kingdom"
$ws.Range("E14").Value = "Flambards Theme Park, Clodgey Ln, Helston TR13 0QA"
$ws.Range("E15").Value = "Flamingo Land Resort, Kirby Misperton, Malton"
$ws.Range("E16").Value = "Address not found"
$ws.Range("E17").Value = "Great Yarmouth Pleasure Beach, The Pleasure Beach, S Beach Parade, Great Yarmouth NR30 3EH, United Kingdom"
$ws.Range("E18").Value = "GreenWood Family Park, Bush Rd, Y Felinheli LL56 4QN"
$ws.Range("E19").Value = "Gulliver's Kingdom, Temple Walk, Matlock Bath DE4 3PG, United Kingdom"
$ws.Range("E20").Value = "Address not found"
$ws.Range("E21").Value = "Gulliver's Valley Theme Park, Mansfield Rd, Sheffield S26 5QW"
$ws.Range("E22").Value = "Gulliver's World Theme Park, Shackleton Cl, Old Hall, Warrington WA5 9YZ, United Kingdom"
$ws.Range("E23").Value = "Harbour Park Amusements, Windmill Rd, Littlehampton BN17 5LH"
$ws.Range("E24").Value = "Address not found"
$ws.Range("E25").Value = "Joyland, Marine Parade, Great Yarmouth NR30 2DL"
$ws.Range("E26").Value = "Address not found"
$ws.Range("E27").Value = "Lightwater Valley Family Adventure Park, Water Ln, North Stainley, Ripon HG4 3HT, United Kingdom"
$ws.Range("E28").Value = "Address not found"
$ws.Range("E29").Value = "Oakwood Theme Park, Canaston Bridge, Narberth SA67 8DE"
$ws.Range("E30").Value = "Ocean Beach Pleasure Park, Sea Rd, South Shields NE33 2LD, United Kingdom"
$ws.Range("E31").Value = "Paultons Park Home of Peppa Pig World, Romsey SO51 6AL, United Kingdom"
$ws.Range("E32").Value = "Pettitts Animal Adventure Park, Church Rd, Reedham, Norwich NR13 3UA"
$ws.Range("E33").Value = "Pleasurewood Hills Family Theme Park, Leisure Way, Lowestoft NR32 5DZ"
$ws.Range("E34").Value = "ROARR!, Lenwade, Norwich NR9 5JE"
$ws.Range("E35").Value = "Southport Pleasureland, Marine Dr, Southport PR8 1RX, United Kingdom"
$ws.Range("E36").Value = "Sundown Adventureland, Treswell Rd, Rampton, Retford DN22 0HX"
$ws.Range("E37").Value = "The BIG Sheep, Abbotsham Rd, Abbotsham, Bideford EX39 5AP"
$ws.Range("E38").Value = "The Milky Way Adventure Park, Bideford EX39 5RY, United Kingdom"
$ws.Range("E39").Value = "Address not found"
$ws.Range("E40").Value = "Address not found"
$ws.Range("E41").Value = "West Midlands Safari Park, Spring Grove, Bewdley DY12 1LF, United Kingdom"
$ws.Range("E42").Value = "Robin Hoods Wheelgate Park, Mansfield Rd, Farnsfield, Newark NG22 8HX, United Kingdom"
$ws.Range("E43").Value = "Wicksteed Park, Barton Rd, Kettering NN15 6NJ, United Kingdom"
$ws.Range("E44").Value = "Woodlands Family Theme Park, Woodlands Leisure Park A3122, Blackawton, Totnes TQ9 7DQ"
